# GanttChart.xlsx update: renamed the project title and several tasks,
# re-scoped a few task durations (the End Date / Range-for-Gantt-Chart
# cells recalc automatically from these), and left the selection where
# the editing session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task names (B6:B11) --------------------------------------------------
$ws.Range("B8").Value  = "Develop Service Layer"
$ws.Range("B6").Value  = "Develop Repository Layer"
$ws.Range("B7").Value  = "Develop nHibernate Mappings"
$ws.Range("B9").Value  = "Develop Controllers "
$ws.Range("B10").Value = "Develop Views"
$ws.Range("B11").Value = "Site CSS Design "

# --- Project title ---------------------------------------------------------
$ws.Range("B2").Value = "TravelMe"

# --- Task durations (D5:D11) ---------------------------------------------
# End dates / the Range-for-Gantt-Chart H5 cell are formula driven and
# recalculate automatically from these.
$ws.Range("D5").Value  = 3
$ws.Range("D7").Value  = 5
$ws.Range("D8").Value  = 7
$ws.Range("D10").Value = 7
$ws.Range("D11").Value = 14

# --- Leftover selection from editing the chart/table ----------------------
$ws.Range("H11").Select()

